$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 974
$ws.Range("C3").Value = 1889
$ws.Range("D3").Value = 4404
$ws.Range("E3").Value = 7656
$ws.Range("F3").Value = 9277
$ws.Range("G3").Value = 15900
$ws.Range("B8").Value = 8626
$ws.Range("C8").Value = 14300
$ws.Range("D8").Value = 32200
$ws.Range("E8").Value = 54900
$ws.Range("F8").Value = 95800
$ws.Range("G8").Value = 120000
$ws.Range("B13").Value = 7236
$ws.Range("C13").Value = 13200
$ws.Range("D13").Value = 13200
$ws.Range("E13").Value = 17300
$ws.Range("F13").Value = 17500
$ws.Range("G13").Value = 14400
$ws.Range("B18").Value = 143000
$ws.Range("C18").Value = 167000
$ws.Range("D18").Value = 323000
$ws.Range("E18").Value = 416000
$ws.Range("F18").Value = 384000
$ws.Range("G18").Value = 216000
$ws.Range("B23").Value = 18100
$ws.Range("C23").Value = 4496
$ws.Range("D23").Value = 5764
$ws.Range("E23").Value = 7571
$ws.Range("F23").Value = 11200
$ws.Range("G23").Value = 6858
$ws.Range("B28").Value = 136000
$ws.Range("C28").Value = 153000
$ws.Range("D28").Value = 247000
$ws.Range("E28").Value = 188000
$ws.Range("F28").Value = 616000
$ws.Range("G28").Value = 237000
$ws.Range("B33").Value = 7641
$ws.Range("C33").Value = 9416
$ws.Range("D33").Value = 11900
$ws.Range("E33").Value = 13300
$ws.Range("F33").Value = 13400
$ws.Range("G33").Value = 13100
$ws.Range("B38").Value = 132000
$ws.Range("C38").Value = 187000
$ws.Range("D38").Value = 297000
$ws.Range("E38").Value = 354000
$ws.Range("F38").Value = 366000
$ws.Range("G38").Value = 296000
